$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "38.789.80"
$rng.ClearFormats()
$rng = $ws.Range("E2")
$rng.NumberFormat = "@"
$rng.Value = "  +1.19%  "
$rng.ClearFormats()
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "2.101.01"
$rng.ClearFormats()
$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = "  -0.07%  "
$rng.ClearFormats()
$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$rng.Value = "  -0.04%  "
$rng.ClearFormats()
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "228.75"
$rng.ClearFormats()
$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$rng.Value = "  -0.17%  "
$rng.ClearFormats()
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "62.36"
$rng.ClearFormats()
$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$rng.Value = "  +1.87%  "
$rng.ClearFormats()
$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$rng.Value = "  -0.08%  "
$rng.ClearFormats()
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "0.387"
$rng.ClearFormats()
$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = "  +1.69%  "
$rng.ClearFormats()
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "0.0841"
$rng.ClearFormats()
$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = "  -0.63%  "
$rng.ClearFormats()
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "0.104"
$rng.ClearFormats()
$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$rng.Value = "  -0.07%  "
$rng.ClearFormats()
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "15.83"
$rng.ClearFormats()
$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$rng.Value = "  +7.26%  "
$rng.ClearFormats()
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "2.413.42"
$rng.ClearFormats()
$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$rng.Value = "  -0.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "22.20"
$rng.ClearFormats()
$rng = $ws.Range("E14")
$rng.NumberFormat = "@"
$rng.Value = "  -0.77%  "
$rng.ClearFormats()
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "0.806"
$rng.ClearFormats()
$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$rng.Value = "  +3.53%  "
$rng.ClearFormats()
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "5.51"
$rng.ClearFormats()
$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$rng.Value = "  +0.30%  "
$rng.ClearFormats()
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "2.113.09"
$rng.ClearFormats()
$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$rng.Value = "  +1.71%  "
$rng.ClearFormats()
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "38.798.85"
$rng.ClearFormats()
$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$rng.Value = "  +1.39%  "
$rng.ClearFormats()
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "72.07"
$rng.ClearFormats()
$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = "  +2.31%  "
$rng.ClearFormats()
$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$rng.Value = "  +0.27%  "
$rng.ClearFormats()
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "0.0₃0839"
$rng.ClearFormats()
$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = "  +0.56%  "
$rng.ClearFormats()
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "227.95"
$rng.ClearFormats()
$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$rng.Value = "  +1.48%  "
$rng.ClearFormats()
$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = "  +0.01%  "
$rng.ClearFormats()
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "2.37"
$rng.ClearFormats()
$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$rng.Value = "  -2.95%  "
$rng.ClearFormats()
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "2.33"
$rng.ClearFormats()
$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$rng.Value = "  +0.45%  "
$rng.ClearFormats()
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "172.00"
$rng.ClearFormats()
$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$rng.Value = "  +1.26%  "
$rng.ClearFormats()
$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$rng.Value = "  +1.61%  "
$rng.ClearFormats()
$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = "  +6.64%  "
$rng.ClearFormats()
$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = "  +4.05%  "
$rng.ClearFormats()
$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$rng.Value = "  +1.67%  "
$rng.ClearFormats()
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "2.47"
$rng.ClearFormats()
$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$rng.Value = "  +3.39%  "
$rng.ClearFormats()
$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$rng.Value = "  +0.82%  "
$rng.ClearFormats()
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "4.53"
$rng.ClearFormats()
$rng = $ws.Range("E33")
$rng.NumberFormat = "@"
$rng.Value = "  +2.06%  "
$rng.ClearFormats()
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "4.77"
$rng.ClearFormats()
$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$rng.Value = "  +1.20%  "
$rng.ClearFormats()
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "0.0620"
$rng.ClearFormats()
$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$rng.Value = "  +2.38%  "
$rng.ClearFormats()
$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$rng.Value = "  +2.58%  "
$rng.ClearFormats()
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "2.41"
$rng.ClearFormats()
$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$rng.Value = "  +0.96%  "
$rng.ClearFormats()
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "3.55"
$rng.ClearFormats()
$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$rng.Value = "  +0.64%  "
$rng.ClearFormats()
$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$rng.Value = "  -0.06%  "
$rng.ClearFormats()
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "18.38"
$rng.ClearFormats()
$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$rng.Value = "  +1.51%  "
$rng.ClearFormats()
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "0.0228"
$rng.ClearFormats()
$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$rng.Value = "  +4.04%  "
$rng.ClearFormats()
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "102.13"
$rng.ClearFormats()
$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = "  +2.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "1.533.07"
$rng.ClearFormats()
$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$rng.Value = "  -1.01%  "
$rng.ClearFormats()
$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$rng.Value = "  -0.96%  "
$rng.ClearFormats()
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "7.78"
$rng.ClearFormats()
$rng = $ws.Range("E45")
$rng.NumberFormat = "@"
$rng.Value = "  +3.83%  "
$rng.ClearFormats()
$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$rng.Value = "  +0.17%  "
$rng.ClearFormats()
$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = "  +1.61%  "
$rng.ClearFormats()
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "4.15"
$rng.ClearFormats()
$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$rng.Value = "  -0.10%  "
$rng.ClearFormats()
$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$rng.Value = "  +1.28%  "
$rng.ClearFormats()
$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$rng.Value = "  -1.24%  "
$rng.ClearFormats()
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "2.296.53"
$rng.ClearFormats()
$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$rng.Value = "  -0.21%  "
$rng.ClearFormats()
